# Atualização 30/06/2020 - 01/07/2020
# Adds two new rows (110 and 111) of data to the isolation-index sheet,
# continuing the existing table which ends at row 109 (6/28/2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110: 6/29/2020, 46%, 0.46, Segunda-feira
# Leading "'" forces Date/Percent-looking text to stay literal text
# (matching the existing rows, which are all stored as plain text),
# rather than being auto-converted to a date serial / percent number.
$ws.Range("A110").Value = "'6/29/2020"
$ws.Range("B110").Value = "'46%"
$ws.Range("C110").Value = 0.46
$ws.Range("D110").Value = "Segunda-feira"

# Row 111: 6/30/2020, 46%, 0.46, Terça-feira
$ws.Range("A111").Value = "'6/30/2020"
$ws.Range("B111").Value = "'46%"
$ws.Range("C111").Value = 0.46
$ws.Range("D111").Value = "Terça-feira"

# Drop the implicit "quote prefix" formatting created by the leading
# apostrophe above so the new cells keep the same (default/no-style)
# formatting as every other data row in the sheet.
$ws.Range("A110:D111").ClearFormats()
